$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (td_sim_1) and column D (record_atd) values as corrected
# in the "Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)" fix.

$ws.Range("C2").Value = 753
$ws.Range("D2").Value = 720.5

$ws.Range("C3").Value = 742
$ws.Range("D3").Value = 739

$ws.Range("C4").Value = 1061
$ws.Range("D4").Value = 987

$ws.Range("C5").Value = 156
$ws.Range("D5").Value = 90

$ws.Range("C6").Value = 482
$ws.Range("D6").Value = 351.5

$ws.Range("C7").Value = 22
$ws.Range("D7").Value = 504.5

$ws.Range("C8").Value = 132
$ws.Range("D8").Value = 132

$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 502.5

$ws.Range("C10").Value = 161
$ws.Range("D10").Value = 95

$ws.Range("C11").Value = 392.1111111111111
